# Generate Report for Handback
# Adds a new handback entry (c500e649-6723-4740-b5a5-7740bf5d0c94) as row 3
# on the Overview, zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$guidNew = "c500e649-6723-4740-b5a5-7740bf5d0c94"
$shaZhCn = "a09e22d16c9bdd0485fe61508e945aff51a3c410"
$shaDeDe = "a09e22d16c9bdd0485fe61508e945aff51a3c410"

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.ListRows.Add() | Out-Null

$wsOverview.Range("A3").Value = "$guidNew.md"
$wsOverview.Range("B3").Value = "e2e\$guidNew.md"
$wsOverview.Range("C3").Value = ".md"
$wsOverview.Range("E3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("G3").Value = "2016-09-07 08:38:16"

$wsOverview.Range("B3").Font.Color = 15570276
$wsOverview.Range("B3").Font.Underline = $true
$wsOverview.Range("G3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2f72b01d371afdcfdc8c7b384f29b2f0c95db82e/e2e/$guidNew.md", "", "", "e2e\$guidNew.md") | Out-Null

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$loZhCn = $wsZhCn.ListObjects.Item(1)
$loZhCn.ListRows.Add() | Out-Null

$wsZhCn.Range("A3").Value = "$guidNew.md"
$wsZhCn.Range("B3").Value = ".md"
$wsZhCn.Range("C3").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("D3").Value = "e2e"
$wsZhCn.Range("E3").Value = "ht"
$wsZhCn.Range("F3").Value = "'True"
$wsZhCn.Range("G3").Value = "$guidNew.$shaZhCn.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-09-07 08:37:56"
$wsZhCn.Range("I3").Value = "$guidNew.md"
$wsZhCn.Range("J3").Value = "$guidNew.$shaZhCn.zh-cn.xlf"
$wsZhCn.Range("K3").Value = "2016-09-07 08:38:53"
$wsZhCn.Range("L3").Value = "'"
$wsZhCn.Range("M3").Value = "'True"
$wsZhCn.Range("N3").Value = "'"
$wsZhCn.Range("O3").Value = "'False"
$wsZhCn.Range("P3").Value = "'"

$wsZhCn.Range("A3").Font.Color = 15570276
$wsZhCn.Range("A3").Font.Underline = $true
$wsZhCn.Range("I3").Font.Color = 15570276
$wsZhCn.Range("I3").Font.Underline = $true
$wsZhCn.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2f72b01d371afdcfdc8c7b384f29b2f0c95db82e/e2e/$guidNew.md", "", "", "$guidNew.md") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/8c211647b308c9bf35320950983744e1cfa03adc/e2e/$guidNew.md", "", "", "$guidNew.md") | Out-Null

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$loDeDe = $wsDeDe.ListObjects.Item(1)
$loDeDe.ListRows.Add() | Out-Null

$wsDeDe.Range("A3").Value = "$guidNew.md"
$wsDeDe.Range("B3").Value = ".md"
$wsDeDe.Range("C3").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("D3").Value = "e2e"
$wsDeDe.Range("E3").Value = "ht"
$wsDeDe.Range("F3").Value = "'True"
$wsDeDe.Range("G3").Value = "$guidNew.$shaDeDe.de-de.xlf"
$wsDeDe.Range("H3").Value = "2016-09-07 08:38:16"
$wsDeDe.Range("I3").Value = "$guidNew.md"
$wsDeDe.Range("J3").Value = "$guidNew.$shaDeDe.de-de.xlf"
$wsDeDe.Range("K3").Value = "2016-09-07 08:39:20"
$wsDeDe.Range("L3").Value = "'"
$wsDeDe.Range("M3").Value = "'True"
$wsDeDe.Range("N3").Value = "'"
$wsDeDe.Range("O3").Value = "'False"
$wsDeDe.Range("P3").Value = "'"

$wsDeDe.Range("A3").Font.Color = 15570276
$wsDeDe.Range("A3").Font.Underline = $true
$wsDeDe.Range("I3").Font.Color = 15570276
$wsDeDe.Range("I3").Font.Underline = $true
$wsDeDe.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2f72b01d371afdcfdc8c7b384f29b2f0c95db82e/e2e/$guidNew.md", "", "", "$guidNew.md") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/28febd96951ee9922bbc33a6ee84192ab5eed13b/e2e/$guidNew.md", "", "", "$guidNew.md") | Out-Null
